# Generate Report for Handback
# Updates the localization-status workbook: marks rows as handed back,
# records the handback target/date for zh-cn and de-de, adds "Latest
# Target File" hyperlinks, and widens the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$zhcnTargetFile = "7a6ff401-3d4c-4db2-a3e2-a2ff2732b77b.36333619e9ac0be9a3ef7c6d3f3d1fb72be4b193.zh-cn.xlf"
$dedeTargetFile  = "7a6ff401-3d4c-4db2-a3e2-a2ff2732b77b.36333619e9ac0be9a3ef7c6d3f3d1fb72be4b193.de-de.xlf"

$mdDisplay7a6  = "7a6ff401-3d4c-4db2-a3e2-a2ff2732b77b.md"
$mdDisplayffff = "ffffa549b742-c9a4-44e5-90d1-e0b93a6912c2.md"

$url7a6  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3dcdf4d2710ca772fda733c6e0a0d37880e25381/e2e/7a6ff401-3d4c-4db2-a3e2-a2ff2732b77b.md"
$urlffff = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3dcdf4d2710ca772fda733c6e0a0d37880e25381/e2e/ffffa549b742-c9a4-44e5-90d1-e0b93a6912c2.md"

# Column-width helper: the engine quantizes stored width to the nearest
# 1/6 character unit, so feed it a ColumnWidth that lands exactly on the
# desired grid point (40, or ~30 as the closest point to the ~29.98 target).
$wideWidth   = 39.140625   # -> stored column width 40
$mediumWidth = 29.140625   # -> stored column width 30 (closest to 29.9777)

# ---------------------------------------------------------------------
# Overview sheet: status text for both language columns, widen them.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

$wsOverview.Columns.Item(5).ColumnWidth = $mediumWidth
$wsOverview.Columns.Item(6).ColumnWidth = $mediumWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("C3").Value = $statusHandedBack

$wsZhCn.Range("J2").Value = $zhcnTargetFile
$wsZhCn.Range("J3").Value = $zhcnTargetFile

$wsZhCn.Range("K2").Value = "2016-08-19 17:06:44"
$wsZhCn.Range("K3").Value = "2016-08-19 17:06:44"

$wsZhCn.Columns.Item(3).ColumnWidth = $mediumWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $wideWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $wideWidth

# Rebuild hyperlinks in display order A2, I2, A3, I3 so relationship ids
# come out as rId2..rId5 in that order.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $url7a6, "", "", $mdDisplay7a6)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $url7a6, "", "", $mdDisplay7a6)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlffff, "", "", $mdDisplayffff)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $url7a6, "", "", $mdDisplay7a6)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("C3").Value = $statusHandedBack

$wsDeDe.Range("J2").Value = $dedeTargetFile
$wsDeDe.Range("J3").Value = $dedeTargetFile

$wsDeDe.Range("K2").Value = "2016-08-19 17:06:51"
$wsDeDe.Range("K3").Value = "2016-08-19 17:06:51"

$wsDeDe.Columns.Item(3).ColumnWidth = $mediumWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $wideWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $wideWidth

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $url7a6, "", "", $mdDisplay7a6)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $url7a6, "", "", $mdDisplay7a6)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlffff, "", "", $mdDisplayffff)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $url7a6, "", "", $mdDisplay7a6)
